$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells we touch to Text format so that numeric-looking
# strings (e.g. "569.94") are preserved exactly as text, matching the source
# workbook which stores these as inline strings, not numbers.
$ws.Range('D2:D6').NumberFormat = '@'
$ws.Range('D8:D9').NumberFormat = '@'
$ws.Range('D11:D22').NumberFormat = '@'
$ws.Range('D25:D31').NumberFormat = '@'
$ws.Range('D34:D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45:D51').NumberFormat = '@'

$ws.Range('D2').Value = '60.536.03'
$ws.Range('E2').Value = '  -2.56%  '

$ws.Range('D3').Value = '2.408.31'
$ws.Range('E3').Value = '  -1.16%  '

$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.39%  '

$ws.Range('D5').Value = '569.94'
$ws.Range('E5').Value = '  -1.61%  '

$ws.Range('D6').Value = '138.83'
$ws.Range('E6').Value = '  -3.66%  '

$ws.Range('E7').Value = '  +0.30%  '

$ws.Range('D8').Value = '0.524'
$ws.Range('E8').Value = '  -1.37%  '

$ws.Range('D9').Value = '2.390.33'
$ws.Range('E9').Value = '  -1.94%  '

$ws.Range('E10').Value = '  +0.14%  '

$ws.Range('D11').Value = '0.160'
$ws.Range('E11').Value = '  +0.09%  '

$ws.Range('D12').Value = '5.06'
$ws.Range('E12').Value = '  -2.85%  '

$ws.Range('D13').Value = '0.339'
$ws.Range('E13').Value = '  -1.59%  '

$ws.Range('D14').Value = '25.76'
$ws.Range('E14').Value = '  -2.31%  '

$ws.Range('D15').Value = '0.0000170'
$ws.Range('E15').Value = '  -1.96%  '

$ws.Range('D16').Value = '2.760.66'
$ws.Range('E16').Value = '  -3.96%  '

$ws.Range('D17').Value = '60.505.32'
$ws.Range('E17').Value = '  -2.54%  '

$ws.Range('D18').Value = '2.378.22'
$ws.Range('E18').Value = '  -2.43%  '

$ws.Range('D19').Value = '10.56'
$ws.Range('E19').Value = '  -2.53%  '

$ws.Range('D20').Value = '7.30'
$ws.Range('E20').Value = '  +1.79%  '

$ws.Range('D21').Value = '320.36'
$ws.Range('E21').Value = '  -2.36%  '

$ws.Range('D22').Value = '4.02'
$ws.Range('E22').Value = '  -2.27%  '

$ws.Range('E23').Value = '  +1.08%  '

$ws.Range('E24').Value = '  +0.33%  '

$ws.Range('D25').Value = '1.87'
$ws.Range('E25').Value = '  -5.13%  '

$ws.Range('D26').Value = '64.53'
$ws.Range('E26').Value = '  -1.73%  '

$ws.Range('D27').Value = '8.45'
$ws.Range('E27').Value = '  -9.49%  '

$ws.Range('D28').Value = '574.07'
$ws.Range('E28').Value = '  -5.91%  '

$ws.Range('D29').Value = '2.495.48'
$ws.Range('E29').Value = '  -2.17%  '

$ws.Range('D30').Value = '0.0₃0920'
$ws.Range('E30').Value = '  -3.53%  '

$ws.Range('D31').Value = '7.87'
$ws.Range('E31').Value = '  -1.46%  '

$ws.Range('E32').Value = '  -6.11%  '

$ws.Range('E33').Value = '  -2.95%  '

$ws.Range('D34').Value = '0.133'
$ws.Range('E34').Value = '  -4.84%  '

$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.17%  '

$ws.Range('D36').Value = '4.60'
$ws.Range('E36').Value = '  -6.53%  '

$ws.Range('D37').Value = '0.368'
$ws.Range('E37').Value = '  -2.24%  '

$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = '149.21'
$ws.Range('E38').Value = '  -0.07%  '

$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '1.38'
$ws.Range('E39').Value = '  -4.85%  '

$ws.Range('D40').Value = '18.18'
$ws.Range('E40').Value = '  -1.19%  '

$ws.Range('D41').Value = '5.11'
$ws.Range('E41').Value = '  -4.78%  '

$ws.Range('E42').Value = '  +0.07%  '

$ws.Range('D43').Value = '1.66'
$ws.Range('E43').Value = '  -4.68%  '

$ws.Range('E44').Value = '  -4.16%  '

$ws.Range('D45').Value = '2.34'
$ws.Range('E45').Value = '  -5.32%  '

$ws.Range('D46').Value = '0.0₆0290'
$ws.Range('E46').Value = '  +16.27%  '

$ws.Range('D47').Value = '140.59'
$ws.Range('E47').Value = '  -1.37%  '

$ws.Range('D48').Value = '3.51'
$ws.Range('E48').Value = '  -3.32%  '

$ws.Range('D49').Value = '0.585'
$ws.Range('E49').Value = '  -3.33%  '

$ws.Range('D50').Value = '0.0502'
$ws.Range('E50').Value = '  -4.03%  '

$ws.Range('D51').Value = '19.23'
$ws.Range('E51').Value = '  -1.85%  '
